$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.527.96"
$ws.Range("E2").Value = "  +0.58%  "
$ws.Range("D3").Value = "2.485.38"
$ws.Range("E3").Value = "  +0.95%  "
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").Value = "'312.91"
$ws.Range("E5").Value = "  +0.55%  "
$ws.Range("D6").Value = "'93.22"
$ws.Range("E6").Value = "  -1.00%  "
$ws.Range("D7").Value = "'0.547"
$ws.Range("E7").Value = "  -0.74%  "
$ws.Range("E8").Value = "  -0.17%  "
$ws.Range("E9").Value = "  -0.66%  "
$ws.Range("D10").Value = "'32.57"
$ws.Range("E10").Value = "  -2.47%  "
$ws.Range("D11").Value = "'0.0784"
$ws.Range("E11").Value = "  +0.56%  "
$ws.Range("E12").Value = "  +2.06%  "
$ws.Range("D13").Value = "2.866.12"
$ws.Range("E13").Value = "  +1.00%  "
$ws.Range("E14").Value = "  -1.40%  "
$ws.Range("D15").Value = "'15.50"
$ws.Range("E15").Value = "  +6.15%  "
$ws.Range("D16").Value = "2.497.13"
$ws.Range("E16").Value = "  +2.23%  "
$ws.Range("D17").Value = "'0.753"
$ws.Range("E17").Value = "  -4.05%  "
$ws.Range("D18").Value = "41.558.24"
$ws.Range("E18").Value = "  +0.79%  "
$ws.Range("D19").Value = "'6.34"
$ws.Range("E19").Value = "  +0.34%  "
$ws.Range("D20").Value = "0.0₃0931"
$ws.Range("E20").Value = "  +1.50%  "
$ws.Range("E21").Value = "  +4.57%  "
$ws.Range("D22").Value = "'11.20"
$ws.Range("E22").Value = "  -1.99%  "
$ws.Range("D23").Value = "'235.85"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("E24").Value = "  -2.55%  "
$ws.Range("E25").Value = "  -0.15%  "
$ws.Range("E26").Value = "  -0.99%  "
$ws.Range("D27").Value = "'24.79"
$ws.Range("E27").Value = "  +1.83%  "
$ws.Range("D28").Value = "'2.23"
$ws.Range("E28").Value = "  +0.22%  "
$ws.Range("D29").Value = "'9.65"
$ws.Range("E29").Value = "  -0.33%  "
$ws.Range("D30").Value = "'36.31"
$ws.Range("E30").Value = "  +0.70%  "
$ws.Range("D31").Value = "'157.08"
$ws.Range("E31").Value = "  +2.76%  "
$ws.Range("D32").Value = "'5.44"
$ws.Range("E32").Value = "  -1.81%  "
$ws.Range("B33").Value = "Celestia"
$ws.Range("C33").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D33").Value = "'18.20"
$ws.Range("E33").Value = "  +7.04%  "
$ws.Range("B34").Value = "WEMIXToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").Value = "'2.57"
$ws.Range("E34").Value = "  -0.75%  "
$ws.Range("E35").Value = "  +0.18%  "
$ws.Range("E36").Value = "  -5.42%  "
$ws.Range("E37").Value = "  -1.60%  "
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "'0.105"
$ws.Range("E38").Value = "  +2.85%  "
$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").Value = "'1.84"
$ws.Range("E39").Value = "  -2.77%  "
$ws.Range("E40").Value = "  -0.05%  "
$ws.Range("D41").Value = "'4.12"
$ws.Range("E41").Value = "  -2.36%  "
$ws.Range("E42").Value = "  -0.15%  "
$ws.Range("D43").Value = "'19.86"
$ws.Range("E43").Value = "  -4.46%  "
$ws.Range("D44").Value = "1.962.26"
$ws.Range("E44").Value = "  +0.19%  "
$ws.Range("E45").Value = "  +0.07%  "
$ws.Range("D46").Value = "'2.96"
$ws.Range("E46").Value = "  -3.19%  "
$ws.Range("D47").Value = "'8.84"
$ws.Range("E47").Value = "  +2.35%  "
$ws.Range("D48").Value = "2.726.82"
$ws.Range("E48").Value = "  +0.90%  "
$ws.Range("D49").Value = "'96.23"
$ws.Range("D50").Value = "'67.53"
$ws.Range("E50").Value = "  -3.20%  "
$ws.Range("D51").Value = "'73.43"
$ws.Range("E51").Value = "  -2.97%  "
